# Apply updates described by the commit:
# "got supporting entity visualization now differentiates between
#  non-numerical and numerical answers."

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "SemScores Analysis"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("SemScores Analysis")

# Row 24: was boolean FALSE -> becomes text "Non-numerical"; refresh stats
$ws1.Range("B24").Value = "Non-numerical"
$ws1.Range("C24").Value = 0.3247443805396292
$ws1.Range("D24").Value = 0.03750033333597924

# Row 25: was boolean TRUE -> becomes text "Got_supporting_entities"
$ws1.Range("B25").Value = "Got_supporting_entities"

# New row 26: No_supporting_entities stats
$ws1.Range("A26").Value = "got_supporting_ents"
$ws1.Range("B26").Value = "No_supporting_entities"
$ws1.Range("C26").Value = 0.224317392432801
$ws1.Range("D26").Value = 0.02082416114565826

# New row 27: overall stats
$ws1.Range("A27").Value = "Overall"
$ws1.Range("B27").Value = "ALL"
$ws1.Range("C27").Value = 0.32951964310579
$ws1.Range("D27").Value = 0.03832279077079733

# ---------------------------------------------------------------------
# Sheet 2: "Best Performers" - column R holds got_supporting_ents flags
# True  -> Got_supporting_entities
# False -> Non-numerical
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Best Performers")
$sheet2Values = @{
    2  = "Got_supporting_entities"
    3  = "Got_supporting_entities"
    4  = "Got_supporting_entities"
    5  = "Got_supporting_entities"
    6  = "Non-numerical"
    7  = "Got_supporting_entities"
    8  = "Non-numerical"
    9  = "Non-numerical"
    10 = "Got_supporting_entities"
    11 = "Got_supporting_entities"
}
foreach ($r in $sheet2Values.Keys) {
    $ws2.Range("R$r").Value = $sheet2Values[$r]
}

# ---------------------------------------------------------------------
# Sheet 3: "Worst Performers" - column R holds got_supporting_ents flags
# False -> Non-numerical, except row 6 -> No_supporting_entities
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Worst Performers")
for ($r = 2; $r -le 11; $r++) {
    if ($r -eq 6) {
        $ws3.Range("R$r").Value = "No_supporting_entities"
    } else {
        $ws3.Range("R$r").Value = "Non-numerical"
    }
}
